$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.625.76'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '3.445.67'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.12'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.00'
$ws.Range("E6").Value = '  +2.60%  '
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.688'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.48'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.48'
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.98'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '3.341.06'
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("D16").Value = '62.519.50'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.43'
$ws.Range("E17").Value = '  +5.62%  '
$ws.Range("E18").Value = '  -1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000133'
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("E20").Value = '  -4.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '84.21'
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '315.22'
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.97'
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.74'
$ws.Range("E25").Value = '  +8.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.85'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.31'
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.75'
$ws.Range("E28").Value = '  +3.21%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.174'
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '42.38'
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  -3.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0487'
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.58'
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.997'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("E38").Value = '  -3.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.319'
$ws.Range("E40").Value = '  +12.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.00'
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.125'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.50'
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.05'
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.89'
$ws.Range("E45").Value = '  -3.26%  '
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.41'
$ws.Range("E47").Value = '  -3.84%  '
$ws.Range("D48").Value = '2.132.00'
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.32'
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.74'
$ws.Range("E51").Value = '  +22.83%  '
